$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "98.500.34"
Set-TextValue $ws.Range("E2") "  -0.52%  "

Set-TextValue $ws.Range("D3") "3.341.85"
Set-TextValue $ws.Range("E3") "  -1.77%  "

Set-TextValue $ws.Range("E4") "  +0.08%  "

Set-TextValue $ws.Range("D5") "262.94"
Set-TextValue $ws.Range("E5") "  +0.74%  "

Set-TextValue $ws.Range("D6") "645.61"
Set-TextValue $ws.Range("E6") "  +1.45%  "

Set-TextValue $ws.Range("D7") "1.51"
Set-TextValue $ws.Range("E7") "  +7.82%  "

Set-TextValue $ws.Range("D8") "0.458"
Set-TextValue $ws.Range("E8") "  +14.82%  "

Set-TextValue $ws.Range("D9") "1.07"
Set-TextValue $ws.Range("E9") "  +19.82%  "

Set-TextValue $ws.Range("D10") "1.00"
Set-TextValue $ws.Range("E10") "  +0.08%  "

Set-TextValue $ws.Range("D11") "3.339.16"
Set-TextValue $ws.Range("E11") "  -1.73%  "

Set-TextValue $ws.Range("D12") "44.14"
Set-TextValue $ws.Range("E12") "  +20.88%  "

Set-TextValue $ws.Range("D13") "0.207"
Set-TextValue $ws.Range("E13") "  +3.38%  "

Set-TextValue $ws.Range("D14") "0.0000272"
Set-TextValue $ws.Range("E14") "  +8.22%  "

Set-TextValue $ws.Range("D15") "98.298.16"
Set-TextValue $ws.Range("E15") "  -0.34%  "

Set-TextValue $ws.Range("D16") "3.979.73"
Set-TextValue $ws.Range("E16") "  -1.08%  "

Set-TextValue $ws.Range("D17") "5.55"
Set-TextValue $ws.Range("E17") "  -0.66%  "

Set-TextValue $ws.Range("D18") "3.343.75"
Set-TextValue $ws.Range("E18") "  -1.23%  "

Set-TextValue $ws.Range("D19") "7.42"
Set-TextValue $ws.Range("E19") "  +18.56%  "

Set-TextValue $ws.Range("D20") "16.62"
Set-TextValue $ws.Range("E20") "  +8.37%  "

Set-TextValue $ws.Range("D21") "530.62"
Set-TextValue $ws.Range("E21") "  +6.51%  "

Set-TextValue $ws.Range("E22") "  -2.31%  "

Set-TextValue $ws.Range("D23") "10.04"
Set-TextValue $ws.Range("E23") "  +6.07%  "

Set-TextValue $ws.Range("D24") "0.0000212"
Set-TextValue $ws.Range("E24") "  -0.81%  "

Set-TextValue $ws.Range("D25") "0.421"
Set-TextValue $ws.Range("E25") "  +46.59%  "

Set-TextValue $ws.Range("D26") "101.40"
Set-TextValue $ws.Range("E26") "  +13.26%  "

Set-TextValue $ws.Range("D27") "6.05"
Set-TextValue $ws.Range("E27") "  +3.70%  "

Set-TextValue $ws.Range("D28") "12.70"
Set-TextValue $ws.Range("E28") "  +4.25%  "

Set-TextValue $ws.Range("D29") "3.525.57"
Set-TextValue $ws.Range("E29") "  +0.18%  "

Set-TextValue $ws.Range("D30") "0.147"
Set-TextValue $ws.Range("E30") "  +10.13%  "

Set-TextValue $ws.Range("E31") "  +0.11%  "

Set-TextValue $ws.Range("E32") "  +11.37%  "

Set-TextValue $ws.Range("D33") "0.188"
Set-TextValue $ws.Range("E33") "  -3.29%  "

Set-TextValue $ws.Range("E34") "  +0.91%  "

Set-TextValue $ws.Range("D35") "29.03"
Set-TextValue $ws.Range("E35") "  +3.08%  "

Set-TextValue $ws.Range("D36") "0.512"
Set-TextValue $ws.Range("E36") "  +7.55%  "

Set-TextValue $ws.Range("D37") "7.77"
Set-TextValue $ws.Range("E37") "  +4.50%  "

Set-TextValue $ws.Range("E38") "  +2.98%  "

Set-TextValue $ws.Range("E39") "  +2.55%  "

Set-TextValue $ws.Range("D40") "520.93"
Set-TextValue $ws.Range("E40") "  +2.44%  "

Set-TextValue $ws.Range("E41") "  -0.67%  "

$ws.Range("B42").Value = "MantraDAO"
$ws.Range("C42").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
Set-TextValue $ws.Range("D42") "3.88"
Set-TextValue $ws.Range("E42") "  +1.68%  "

$ws.Range("B43").Value = "Fetch.AI"
$ws.Range("C43").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue $ws.Range("D43") "1.31"
Set-TextValue $ws.Range("E43") "  +2.29%  "

Set-TextValue $ws.Range("D44") "3.33"
Set-TextValue $ws.Range("E44") "  -2.39%  "

Set-TextValue $ws.Range("D45") "0.805"
Set-TextValue $ws.Range("E45") "  +1.70%  "

Set-TextValue $ws.Range("E46") "  +0.02%  "

Set-TextValue $ws.Range("D47") "0.0387"
Set-TextValue $ws.Range("E47") "  +18.27%  "

Set-TextValue $ws.Range("D48") "163.92"
Set-TextValue $ws.Range("E48") "  +2.26%  "

Set-TextValue $ws.Range("D49") "2.01"
Set-TextValue $ws.Range("E49") "  +2.43%  "

$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws.Range("D50") "7.67"
Set-TextValue $ws.Range("E50") "  +15.95%  "

$ws.Range("B51").Value = "OKB"
$ws.Range("C51").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws.Range("D51") "49.59"
Set-TextValue $ws.Range("E51") "  +6.10%  "
